$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.207.01"

$ws.Range("D3").Value = "'1.829.60"

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").Value = "'237.37"
$ws.Range("E5").Value = '  -1.20%  '

$ws.Range("D6").Value = "'0.6096"
$ws.Range("E6").Value = '  -4.07%  '

$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = '  +0.17%  '

$ws.Range("D8").Value = "'0.07100"
$ws.Range("E8").Value = '  -5.11%  '

$ws.Range("D9").Value = "'0.2822"
$ws.Range("E9").Value = '  -2.88%  '

$ws.Range("D10").Value = "'23.89"
$ws.Range("E10").Value = '  -4.95%  '

$ws.Range("D11").Value = "'0.07639"
$ws.Range("E11").Value = '  -1.36%  '

$ws.Range("D12").Value = "'1.834.93"
$ws.Range("E12").Value = '  -0.18%  '

$ws.Range("D13").Value = "'4.820"
$ws.Range("E13").Value = '  -3.37%  '

$ws.Range("D14").Value = "'0.6379"
$ws.Range("E14").Value = '  -6.17%  '

$ws.Range("D15").Value = "'0.000009978"
$ws.Range("E15").Value = '  -2.71%  '

$ws.Range("D16").Value = "'2.070.93"
$ws.Range("E16").Value = '  -1.01%  '

$ws.Range("D17").Value = "'79.35"
$ws.Range("E17").Value = '  -3.25%  '

$ws.Range("D18").Value = "'5.960"
$ws.Range("E18").Value = '  -4.80%  '

$ws.Range("D19").Value = "'29.216.32"
$ws.Range("E19").Value = '  -0.80%  '

$ws.Range("D20").Value = "'229.03"
$ws.Range("E20").Value = '  -0.54%  '

$ws.Range("D21").Value = "'11.83"
$ws.Range("E21").Value = '  -4.22%  '

$ws.Range("E22").Value = '  +0.16%  '

$ws.Range("D23").Value = "'7.048"
$ws.Range("E23").Value = '  -5.00%  '

$ws.Range("D24").Value = "'1.003"
$ws.Range("E24").Value = '  +0.29%  '

$ws.Range("D25").Value = "'155.53"
$ws.Range("E25").Value = '  -1.62%  '

$ws.Range("D26").Value = "'8.116"
$ws.Range("E26").Value = '  -4.64%  '

$ws.Range("E27").Value = '  -4.25%  '

$ws.Range("D28").Value = "'16.70"

$ws.Range("D29").Value = "'0.06757"
$ws.Range("E29").Value = '  +3.32%  '

$ws.Range("D30").Value = "'1.485"
$ws.Range("E30").Value = '  +3.87%  '

$ws.Range("E31").Value = '  -2.25%  '

$ws.Range("D32").Value = "'3.857"
$ws.Range("E32").Value = '  -5.35%  '

$ws.Range("D33").Value = "'3.840"
$ws.Range("E33").Value = '  -5.34%  '

$ws.Range("D34").Value = "'1.132"
$ws.Range("E34").Value = '  -0.84%  '

$ws.Range("D35").Value = "'1.742"
$ws.Range("E35").Value = '  -5.41%  '

$ws.Range("D36").Value = "'0.6573"
$ws.Range("E36").Value = '  -6.14%  '

$ws.Range("D37").Value = "'2.556"
$ws.Range("E37").Value = '  -0.81%  '

$ws.Range("D38").Value = "'1.236.93"
$ws.Range("E38").Value = '  -1.17%  '

$ws.Range("D39").Value = "'2.759"
$ws.Range("E39").Value = '  -2.16%  '

$ws.Range("D40").Value = "'0.01765"
$ws.Range("E40").Value = '  -5.10%  '

$ws.Range("D41").Value = "'6.594"
$ws.Range("E41").Value = '  -2.47%  '

$ws.Range("D42").Value = "'0.9250"
$ws.Range("E42").Value = '  -1.44%  '

$ws.Range("E43").Value = '  +0.16%  '

$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = "'101.17"
$ws.Range("E44").Value = '  -0.05%  '

$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").Value = "'1.986.10"
$ws.Range("E45").Value = '  -1.27%  '

$ws.Range("D46").Value = "'63.70"
$ws.Range("E46").Value = '  -2.73%  '

$ws.Range("D47").Value = "'0.00000000116"
$ws.Range("E47").Value = '  -2.22%  '

$ws.Range("D48").Value = "'1.629"
$ws.Range("E48").Value = '  -5.39%  '

$ws.Range("D49").Value = "'8.564"
$ws.Range("E49").Value = '  -4.92%  '

$ws.Range("D50").Value = "'6.541"
$ws.Range("E50").Value = '  -7.54%  '

$ws.Range("E51").Value = '  -5.66%  '
